# The "name" and "code" columns (D = codeforiati:group-name, E = codeforiati:group-code)
# were swapped for every row (header + data): column D now holds the code, column E now
# holds the name. Walk every used row and swap the two cell values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $nameCell = $ws.Cells.Item($r, 4)
    $codeCell = $ws.Cells.Item($r, 5)
    $nameValue = $nameCell.Value2
    $codeValue = $codeCell.Value2
    $nameCell.Value = $codeValue
    $codeCell.Value = $nameValue
}
